# Apply updated naive forecaster QoQ evaluation values
# (ifo GDP component analysis preprocessing) to rows 2-10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  B = -0.0253360835587989;  C = 1.135698984415373;  D = 7.612000487694052;  E = 2.758985409112207;  F = 2.78578550780152;  G = 52 },
    @{ Row = 3;  B = 0.08600705238178873;  C = 1.123115432295725;  D = 5.609313565901017;  E = 2.368398945680608;  F = 2.390387976856151;  G = 51 },
    @{ Row = 4;  B = 0.0363203761572407;   C = 0.9693363125750792; D = 4.390696811956634;  E = 2.095398962478658;  F = 2.116354596420031;  G = 50 },
    @{ Row = 5;  B = 0.1260652424010658;   C = 1.03221972966461;   D = 4.978658809889423;  E = 2.231290839377382;  F = 2.25081259337064;   G = 49 },
    @{ Row = 6;  B = 0.08428106382731648;  C = 1.007929901751503;  D = 5.095203452855081;  E = 2.257255734925726;  F = 2.27955206241149;   G = 48 },
    @{ Row = 7;  B = 0.09619512034022865;  C = 1.131065357075751;  D = 5.749733647272278;  E = 2.397860222630226;  F = 2.42991643159901;   G = 36 },
    @{ Row = 8;  B = 0.1381149580396931;   C = 1.187418322474588;  D = 5.903068426619335;  E = 2.4296231038207;    F = 2.461107737049912;  G = 35 },
    @{ Row = 9;  B = 0.1382699113425551;   C = 1.562389918879535;  D = 9.831370974918348;  E = 3.135501710240061;  F = 3.223266002762909;  G = 18 },
    @{ Row = 10; B = -0.7150897714758779;  C = 1.341686922679743;  D = 6.639370992174167;  E = 2.576697691265735;  F = 2.596308797503721;  G = 11 }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Range("B$r").Value = $u.B
    $ws.Range("C$r").Value = $u.C
    $ws.Range("D$r").Value = $u.D
    $ws.Range("E$r").Value = $u.E
    $ws.Range("F$r").Value = $u.F
    $ws.Range("G$r").Value = $u.G
}
